# Updated cryptos list values (Price column D, Volume(1h) column E).
# Cells whose new Price text reads as a pure number (e.g. "1.00") must be
# forced to Text format first, otherwise Excel auto-coerces the assignment
# to a numeric value and the literal formatting (trailing zero, etc.) is lost.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '66.100.07'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +1.61%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '3.219.46'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +1.33%  '; ForceText = $false }
    @{ Cell = 'D4'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E4'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '602.92'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +4.71%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '151.93'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +1.05%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  +0.01%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '3.217.11'; ForceText = $false }
    @{ Cell = 'E8'; Value = '  +1.49%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.536'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +1.41%  '; ForceText = $false }
    @{ Cell = 'E10'; Value = '  -1.15%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  -0.78%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.511'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +0.99%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '38.60'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +1.41%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '3.748.10'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  +1.46%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '66.135.17'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  +1.51%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '7.42'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  +3.05%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '3.218.74'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +1.74%  '; ForceText = $false }
    @{ Cell = 'E19'; Value = '  +0.56%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '512.97'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -0.16%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '15.84'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +6.24%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '0.738'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  +0.48%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '15.17'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -1.34%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '7.99'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +2.01%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '85.40'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +0.06%  '; ForceText = $false }
    @{ Cell = 'E26'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '9.33'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +2.81%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '3.03'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  +3.98%  '; ForceText = $false }
    @{ Cell = 'E29'; Value = '  +2.07%  '; ForceText = $false }
    @{ Cell = 'E30'; Value = '  +3.96%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '6.83'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +8.05%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '28.20'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +0.24%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  +1.13%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +0.07%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '6.64'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.76%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '55.61'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -0.18%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '0.0924'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +3.15%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '488.69'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +2.08%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.0423'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +0.35%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -2.86%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '8.87'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +2.65%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '3.033.45'; ForceText = $false }
    @{ Cell = 'E42'; Value = '  -2.23%  '; ForceText = $false }
    @{ Cell = 'E43'; Value = '  +0.06%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  +2.26%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.0₃0646'; ForceText = $false }
    @{ Cell = 'E45'; Value = '  +8.33%  '; ForceText = $false }
    @{ Cell = 'E46'; Value = '  +0.45%  '; ForceText = $false }
    @{ Cell = 'E47'; Value = '  -0.82%  '; ForceText = $false }
    @{ Cell = 'E48'; Value = '  +0.08%  '; ForceText = $false }
    @{ Cell = 'E49'; Value = '  +0.53%  '; ForceText = $false }
    @{ Cell = 'E50'; Value = '  +1.58%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '119.48'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -1.59%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Pure-numeric-looking text (e.g. "1.00"): set Text format so Excel
        # keeps it as a literal string instead of coercing to a Double.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        # Restore the default style so no stray formatting is left behind.
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

